# Generate Report for Handoff
#
# The localization pass finished ("In Translation" -> "Ready for handoff")
# and a fresh handoff xliff was produced, so the status + timestamp columns
# are refreshed on the Overview sheet and on each per-locale sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Status: "In Translation" -> "Ready for handoff" ------------------
$ws1.Range("E2").Value2 = "Ready for handoff"   # Overview!zh-cn
$ws1.Range("F2").Value2 = "Ready for handoff"   # Overview!de-de
$ws2.Range("C2").Value2 = "Ready for handoff"   # zh-cn!Status
$ws3.Range("C2").Value2 = "Ready for handoff"   # de-de!Status

# --- Refreshed handoff timestamps --------------------------------------
# Overview!"Latest HO Xliff Generate Date" and de-de!"Latest Handoff Datetime"
# shared the same timestamp before, and still do after the refresh.
$ws1.Range("G2").Value2 = "2016-08-18 04:38:23"
$ws3.Range("H2").Value2 = "2016-08-18 04:38:23"

# zh-cn!"Latest Handoff Datetime" got its own refreshed timestamp.
$ws2.Range("H2").Value2 = "2016-08-18 04:38:18"

# --- Column width grew to fit the longer "Ready for handoff" text ------
$ws1.Columns.Item(5).ColumnWidth = 16.38265482584637   # Overview col E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = 16.38265482584637   # Overview col F (de-de)
$ws2.Columns.Item(3).ColumnWidth = 16.38265482584637   # zh-cn col C (Status)
$ws3.Columns.Item(3).ColumnWidth = 16.38265482584637   # de-de col C (Status)
